# Update the build timestamp embedded in several cells of the workbook.
# Old: "January 30 2026 16.19.47 EST"
# New: "February 02 2026 12.49.33 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

# --- Sheet "About" ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: mines - January 30 (built on " + $newStamp + ")"

$about.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Uskovskaya Coal Mine, Russia, M0864, version 'mines - January 30 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet "Boundaries and methane sources" ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 12; $row++) {
    $cell = $data.Range("S" + $row)
    $current = $cell.Value2
    if ($current -ne $null -and $current.ToString().Contains($oldStamp)) {
        $cell.Value = $current.ToString().Replace($oldStamp, $newStamp)
    }
}
